$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.837.36'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.44%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.857.53'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.89%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.95%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.19%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5077'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3653'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.80%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07176'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.18%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8895'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.39%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.68'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.87%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07531'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.83%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.861.57'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.41%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '91.75'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.43%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.239'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.02%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.13%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008524'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.53%  '

$ws.Range("E18").Value = '  -1.23%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.13%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.884.73'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.41%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.013'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.72%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.094.48'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.23%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.31'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.24%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.445'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.59%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.816'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.88%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '146.05'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.78%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.81'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.49%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.047'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.71%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '112.86'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.641'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.61%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.674'
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09230'
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05095'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.75%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.064'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7336'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.55%  '

$ws.Range("E36").Value = '  -3.03%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.210'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.58%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02008'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.16%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.457'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.64%  '

$ws.Range("E40").Value = '  -1.01%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5308'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.80%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '118.16'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.26%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.498'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.00%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.394'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1474'
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4632'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.71%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9998'
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.884'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.99%  '

$ws.Range("E49").Value = '  -1.29%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.96'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.00%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '62.94'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.47%  '
